# MOSIP TPS Thread setting calculator - RegProc_SyncData
# Performance testing Release 1.3.x v0.5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the worksheet tab ---
$ws.Name = "RegProc_SyncData"

# --- Add reviewer comments (legacy/"Note" comments) left by Deepesh Gurung ---
$excel.UserName = "Deepesh Gurung"

$cmtB6 = $ws.Range("B6").AddComment()
[void]$cmtB6.Text("Deepesh Gurung:" + "`n" + "These come from performance test plan definition ")

$cmtD6 = $ws.Range("D6").AddComment()
[void]$cmtD6.Text("Deepesh Gurung:" + "`n" + "This value is obtained from Jmeter results during previous executions.")

$cmtH6 = $ws.Range("H6").AddComment()
[void]$cmtH6.Text("Deepesh Gurung:" + "`n" + "Apply these values in Jmeter thread setting.")

$cmtK6 = $ws.Range("K6").AddComment()
[void]$cmtK6.Text("Deepesh Gurung:" + "`n" + 'Apply this valie in Jmeter for "RampUp" variable under "user defined variables"')

# --- Adjust the saved view state: scroll position, zoom level and selection ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 4
$win.Zoom = 84
[void]$ws.Range("M10").Select()
